# Adds two missing-item rows (GLAPTIVIA PLUS ... and TAMSULIN 0.4MG 28 CAPS)
# to the "نواقص الأصناف" (missing items) report, keeping the existing
# alphabetical ordering, renumbers the sequence column, refreshes the
# totals cell and bumps the generated-at timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Insert-ReportRow($RowIndex, $Name, $Ratio, $ReqFlag, $Price, $SellPrice, $Trans) {
    $srcRow = $RowIndex - 1

    # Push everything at/after $RowIndex down by one row.
    $ws.Rows("$($RowIndex):$RowIndex").Insert()

    # Clone formatting (number formats, fonts, fills, borders...) from the
    # row immediately above, which carries the same column layout.
    $ws.Range("A$($srcRow):Q$($srcRow)").Copy()
    $ws.Range("A$($RowIndex):Q$($RowIndex)").PasteSpecial(-4122)
    $excel.CutCopyMode = 0

    # Row height matches the surrounding data rows.
    $ws.Rows("$($RowIndex):$RowIndex").RowHeight = 25.5

    # Re-establish the merged regions used by every data row.
    $ws.Range("A$($RowIndex):B$($RowIndex)").Merge()
    $ws.Range("C$($RowIndex):G$($RowIndex)").Merge()
    $ws.Range("H$($RowIndex):K$($RowIndex)").Merge()
    $ws.Range("L$($RowIndex):M$($RowIndex)").Merge()
    $ws.Range("N$($RowIndex):O$($RowIndex)").Merge()

    $ws.Range("C$($RowIndex)").Value = $Name
    $ws.Range("H$($RowIndex)").Value = $Ratio
    $ws.Range("L$($RowIndex)").Value = $ReqFlag
    $ws.Range("N$($RowIndex)").Value = $Price
    $ws.Range("P$($RowIndex)").Value = $SellPrice
    $ws.Range("Q$($RowIndex)").Value = $Trans
}

# "GLAPTIVIA PLUS 50/1000MG 30 F.C. TAB." belongs alphabetically right
# before "HELI-CURE 14 ENTERIC COATED TAB", currently row 21.
Insert-ReportRow 21 "GLAPTIVIA PLUS 50/1000MG 30 F.C. TAB." "0:2" "1" "168.00" "55.4400" "0:1"

# "TAMSULIN 0.4MG 28 CAPS" belongs alphabetically right before
# "VIDROP 2800 I.U./ML ORAL DROPS 15 ML", which (after the insert above)
# now sits on row 39.
Insert-ReportRow 39 "TAMSULIN 0.4MG 28 CAPS" "2:1" "1" "124.00" "62.0000" "0:1"

# Renumber the leading sequence column (A) 1..N for every data row.
$firstDataRow = 7
$lastDataRow = 56
$seq = 1
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $ws.Range("A$($r)").Value = $seq
    $seq++
}

# Refresh the grand-total cell (sum of the "sell price" column) and the
# generated-at timestamp shown in the footer.
$ws.Range("P57").Value = 2303.285
$ws.Range("A58").Value = "Sunday, 22 June, 2025 4:44 PM"
